$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the date columns (C and E) for the data rows so that
# date-looking strings like "2021-05-18" are stored as literal text instead of
# being auto-converted into date serial numbers, matching the source workbook.
$ws.Range("C2:C45").NumberFormat = "@"
$ws.Range("E2:E45").NumberFormat = "@"

$data = @(
    ,@("Wenlambo", " 5 hours", "2021-05-18", "Today", "2021-05-18", 0)
    ,@("SafeMusk", " 10 hours", "2021-05-18", "Today", "2021-05-18", 0)
    ,@("YooShi", " 13 hours", "2021-05-18", "Today", "2021-05-18", 0)
    ,@("Rune", " 15 hours", "2021-05-18", "Today", "2021-05-18", 0)
    ,@("Tiger King", " 15 hours", "2021-05-18", "1 day ", "2021-05-17", 1)
    ,@("GeroWallet", " 17 hours", "2021-05-18", "Today", "2021-05-18", 0)
    ,@("Exohood", " 22 hours", "2021-05-18", "1 day ", "2021-05-17", 1)
    ,@("Relite Finance", "1 day", "2021-05-17", "Today", "2021-05-17", 0)
    ,@("Convex Finance", "1 day", "2021-05-17", "Today", "2021-05-18", 1)
    ,@("420x", "2 days", "2021-05-16", "outdated", "2021-05-07", 9)
    ,@("EthereumMax", "2 days", "2021-05-16", "1 day ", "2021-05-17", 1)
    ,@("Float Protocol: Float", "3 days", "2021-05-15", "1 day ", "2021-05-17", 2)
    ,@("LaunchX", "3 days", "2021-05-14", "outdated", "2021-05-17", 3)
    ,@("Nafter", "4 days", "2021-05-14", "2 days ", "2021-05-15", 1)
    ,@("Alaska Inu", "4 days", "2021-05-14", "outdated", "2021-05-13", 1)
    ,@("Elongate Deluxe", "4 days", "2021-05-14", "outdated", "2021-05-13", 1)
    ,@("SafeZone", "4 days", "2021-05-14", "3 days ", "2021-05-14", 0)
    ,@("Change Our World", "4 days", "2021-05-14", "3 days ", "2021-05-14", 0)
    ,@("Pinknode", "5 days", "2021-05-13", "2 days ", "2021-05-15", 2)
    ,@("Avalaunch", "5 days", "2021-05-13", "3 days ", "2021-05-14", 1)
    ,@("SpaceCorgi", "6 days", "2021-05-11", "outdated", "2021-05-11", 0)
    ,@("PlutoPepe", "7 days", "2021-05-11", "outdated", "2021-05-11", 0)
    ,@("Kabosu", "8 days", "2021-05-10", "outdated", "2021-05-11", 1)
    ,@("MoonTrust", "8 days", "2021-05-10", "outdated", "2021-05-11", 1)
    ,@("Richie", "8 days", "2021-05-10", "outdated", "2021-05-08", 2)
    ,@("Sanshu Inu", "8 days", "2021-05-10", "outdated", "2021-05-11", 1)
    ,@("Blocktyme", "8 days", "2021-05-10", "outdated", "2021-05-09", 1)
    ,@("CateCoin", "8 days", "2021-05-10", "outdated", "2021-05-09", 1)
    ,@("CryptoBlades", "9 days", "2021-05-09", "outdated", "2021-05-09", 0)
    ,@("ShibaPup", "10 days", "2021-05-08", "outdated", "2021-05-11", 3)
    ,@("MoonPirate", "12 days", "2021-05-06", "1 day ", "2021-05-17", 11)
    ,@("Shih Tzu", "12 days", "2021-05-06", "outdated", "2021-05-11", 5)
    ,@("Koduro", "12 days", "2021-05-06", "outdated", "2021-05-15", 9)
    ,@("Ispolink", "12 days", "2021-05-06", "1 day ", "2021-05-17", 11)
    ,@("AlgoPainter", "12 days", "2021-05-06", "outdated", "2021-05-06", 0)
    ,@("PalGold", "13 days", "2021-05-05", "outdated", "2021-05-14", 9)
    ,@("Oviex", "14 days", "2021-05-04", "1 day ", "2021-05-17", 13)
    ,@("Ganesha Token", "15 days", "2021-05-03", "outdated", "2021-05-06", 3)
    ,@("Pussy Financial", "15 days", "2021-05-03", "outdated", "2021-05-07", 4)
    ,@("Gorilla Diamond", "15 days", "2021-05-03", "outdated", "2021-05-14", 11)
    ,@("Lunar Highway", "15 days", "2021-05-03", "outdated", "2021-05-14", 11)
    ,@("CCSwap", "16 days", "2021-05-02", "outdated", "2021-05-08", 6)
    ,@("Altura", "19 days", "2021-04-29", "outdated", "2021-05-07", 8)
    ,@("Token Echo", "20 days", "2021-04-28", "outdated", "2021-05-08", 10)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $row = $row + 1
}

Write-Host "Done, last row written: $($row - 1)"
